# Tarigopula_LabExam03Grading.xlsx - grading rubric update
# - Row 20 (Q12, findNoOfCustomers() method): points-given changed 5 -> 3,
#   and the grading comment swapped to a new, more specific comment.
# - Active selection moved from E37 to G24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update points awarded for question 12 (findNoOfCustomers() method)
$ws.Range("E20").Value = 3

# Update grading comment for that same row with the new rubric text
$ws.Range("F20").Value = "(-7) For incorrect condition for checking if customer exists or not and no need to iterate over products."

# Move the active selection/cursor to G24
$ws.Range("G24").Select()
